$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "TotalShell" column (K) ---

# Header cell K1: same bold/centered/wrap-text look as the other headers
# (copy format from A1) but without the header's bottom border.
$ws.Range("K1").Value = "TotalShell"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").Borders.Item(9).LineStyle = -4142

# Rows 2-13 ("PC" cluster rows): TotalShell = SUM(%_LC, %_SC, %_BC)
$ws.Range("K2").Formula = "=SUM(D2,E2,F2)"
$ws.Range("K3:K13").Formula = "=SUM(D3,E3,F3)"
$ws.Range("K12").Formula = "=SUM(D12,E12,F12)"

# Rows 14-28 ("Image J" cluster rows): hard-coded TotalShell values, styled
# the same way as the existing %_SC values in column E for those rows.
$values = @{
    14 = 64;  15 = 16;  16 = 32;  17 = 36;  18 = 24;
    19 = 100; 20 = 92;  21 = 96;  22 = 60;  23 = 60;
    24 = 88;  25 = 68;  26 = 40;  27 = 100; 28 = 100
}
foreach ($r in 14..28) {
    $ws.Range("E$r").Copy() | Out-Null
    $ws.Range("K$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("K$r").Value = $values[$r]
}

# Widen the new column
$ws.Columns("K").ColumnWidth = 22.7

# Leave the selection where the author left it
$ws.Range("Q12").Select() | Out-Null
